$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update "Förändrad" (changed) date column C for rows 2-6 from 45175 (2023-09-06)
# to 45183 (2023-09-14), matching the original date serial number format.
$ws.Range("C2:C6").Value = 45183
